# Update "Lương" sheet: drop the (stale) column B values entirely, and
# restructure column A to add "Ứng lương ..." rows after each location's
# block plus four new "Tổng lương ..." rows at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Column B no longer carries any values in the new layout - remove it
# outright (this also shrinks the sheet dimension back down to column A).
$ws.Columns.Item(2).Delete()

# Final column-A category list, in row order.
$values = @(
    "Danh mục",
    "Ngày công",
    "Phụ cấp",
    "Lương cơ bản tại CẦN THƠ",
    "Chiết khấu sale chính tại CẦN THƠ",
    "Chiết khấu sale phụ tại CẦN THƠ",
    "Đơn 1 bác sĩ tại CẦN THƠ",
    "Đơn 2 bác sĩ tại CẦN THƠ",
    "Công phụ phẫu 1 tại CẦN THƠ",
    "Công phụ phẫu 2 tại CẦN THƠ",
    "Ứng lương tại CẦN THƠ",
    "Lương cơ bản tại LONG XUYÊN",
    "Chiết khấu sale chính tại LONG XUYÊN",
    "Chiết khấu sale phụ tại LONG XUYÊN",
    "Đơn 1 bác sĩ tại LONG XUYÊN",
    "Đơn 2 bác sĩ tại LONG XUYÊN",
    "Công phụ phẫu 1 tại LONG XUYÊN",
    "Công phụ phẫu 2 tại LONG XUYÊN",
    "Ứng lương tại LONG XUYÊN",
    "Lương cơ bản tại SÓC TRĂNG",
    "Chiết khấu sale chính tại SÓC TRĂNG",
    "Chiết khấu sale phụ tại SÓC TRĂNG",
    "Đơn 1 bác sĩ tại SÓC TRĂNG",
    "Đơn 2 bác sĩ tại SÓC TRĂNG",
    "Công phụ phẫu 1 tại SÓC TRĂNG",
    "Công phụ phẫu 2 tại SÓC TRĂNG",
    "Ứng lương tại SÓC TRĂNG",
    "Tổng lương tại CẦN THƠ",
    "Tổng lương tại LONG XUYÊN",
    "Tổng lương tại SÓC TRĂNG",
    "Tổng lương"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
